$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "November 2015 " -> split "5" into its own run and move the "_GoBack"
#    bookmark so it now sits right after "2015" (before the trailing space).
# ---------------------------------------------------------------------------

# Locate "November 201" so we know exactly where the "5" begins, without
# hard-coding character offsets.
$findRng = $d.Content
$findRng.Find.ClearFormatting()
$findRng.Find.Execute("November 201", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)
$splitPos = $findRng.End              # position right before the "5"

# Nudge a bookmark into existence at that spot and immediately remove it
# again - this forces the engine to keep "November 201" and "5" in separate
# runs (identically-formatted runs are otherwise coalesced back together)
# without leaving any bookmark markup behind.
$d.Bookmarks.Add("TempRunSplit", $d.Range($splitPos, $splitPos)) | Out-Null
$d.Bookmarks.Item("TempRunSplit").Delete()

# The real "_GoBack" bookmark used to live at the very end of the document;
# drop it from there and re-create it right after "...2015" instead.
$d.Bookmarks.Item("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($splitPos + 1, $splitPos + 1)) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Update history" section: add a new "Added IsTrusted support April
#    2018" line, wrapped in a blank paragraph on either side, right after
#    the existing "Initial release May 2016" line.
# ---------------------------------------------------------------------------

# Find the (until now empty, bookmark-only) paragraph that used to carry
# "_GoBack" at the end of the story.
$lastParaIndex = $d.Paragraphs.Count
$targetPara = $d.Paragraphs.Item($lastParaIndex)

# Split a blank paragraph off right before it ...
$splitStart = $targetPara.Range.Start
$d.Range($splitStart, $splitStart).Text = "`r"

# ... the paragraph that used to hold the bookmark is now the next one; put
# the new sentence there (the old bookmark is already gone, see step 1) ...
$textParaIndex = $lastParaIndex + 1
$textPara = $d.Paragraphs.Item($textParaIndex)
$textInsertPos = $textPara.Range.Start
$d.Range($textInsertPos, $textInsertPos).InsertAfter("Added IsTrusted support April 2018")

# ... and finally split another blank paragraph off right after it, just
# before the final section break.
$textPara = $d.Paragraphs.Item($textParaIndex)
$trailingSplitPos = $textPara.Range.End - 1
$d.Range($trailingSplitPos, $trailingSplitPos).Text = "`r"

# ---------------------------------------------------------------------------
# NOTE: the diff also shows two new entries ("Normal Table", "Table Web 1")
# appearing in styles.xml's <w:latentStyles> table. That table is Word's own
# internal "built-in style gallery" bookkeeping - it is refreshed by Word
# itself on save and is not reachable through the Word object model (no
# Application/Document/Styles member exposes w:lsdException), in this
# runtime or in real Word COM automation, so it is intentionally left alone
# here.
# ---------------------------------------------------------------------------
